# Scheduled-runner refresh of cached market-price columns (H:N) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Leve profit" sheets.
# Only raw cached numeric values are touched; no formulas/structure change.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1335.7826
$ws.Range("I40").Value = 1209.1538
$ws.Range("K40").Value = 1209.1538
$ws.Range("M40").Value = -1034.1538
# Row 98
$ws.Range("H98").Value = 5400.696
$ws.Range("I98").Value = 3310.8
$ws.Range("K98").Value = 3310.8
$ws.Range("M98").Value = -1812.8
# Row 104
$ws.Range("H104").Value = 170.5
$ws.Range("I104").Value = 170.5
$ws.Range("K104").Value = 511.5
$ws.Range("M104").Value = 1235.5
# Row 107
$ws.Range("H107").Value = 385.69565
$ws.Range("I107").Value = 402.1579
$ws.Range("J107").Value = 307.5
$ws.Range("K107").Value = 402.1579
$ws.Range("L107").Value = 307.5
$ws.Range("M107").Value = 1517.8421
$ws.Range("N107").Value = -4147.5
# Row 121
$ws.Range("H121").Value = 2163.25
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 2163.25
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 6489.75
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -9983.75
# Row 122
$ws.Range("H122").Value = 5400.696
$ws.Range("I122").Value = 3310.8
$ws.Range("K122").Value = 9932.400000000001
$ws.Range("M122").Value = -7482.400000000001
# Row 141
$ws.Range("H141").Value = 22701.611
$ws.Range("I141").Value = 7730.7144
$ws.Range("J141").Value = 75099.75
$ws.Range("K141").Value = 23192.1432
$ws.Range("L141").Value = 225299.25
$ws.Range("M141").Value = -18012.1432
$ws.Range("N141").Value = -235659.25

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2800
$ws.Range("I61").Value = 2500
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2500
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2288
$ws.Range("N61").Value = -3424
# Row 63
$ws.Range("H63").Value = 3743.1052
$ws.Range("I63").Value = 4319.909
$ws.Range("K63").Value = 4319.909
$ws.Range("M63").Value = -3633.909
# Row 66
$ws.Range("H66").Value = 3743.1052
$ws.Range("I66").Value = 4319.909
$ws.Range("K66").Value = 21599.545
$ws.Range("M66").Value = -18167.545
# Row 110
$ws.Range("H110").Value = 1073.3
$ws.Range("I110").Value = 988.4
$ws.Range("J110").Value = 1328
$ws.Range("K110").Value = 988.4
$ws.Range("L110").Value = 1328
$ws.Range("M110").Value = 1056.6
$ws.Range("N110").Value = -5418
# Row 122
$ws.Range("H122").Value = 1438.5454
$ws.Range("I122").Value = 1536
$ws.Range("K122").Value = 4608
$ws.Range("M122").Value = -2158
# Row 132
$ws.Range("H132").Value = 1528.409
$ws.Range("I132").Value = 979.2778
$ws.Range("K132").Value = 2937.8334
$ws.Range("M132").Value = -407.8334
# Row 136
$ws.Range("H136").Value = 2800
$ws.Range("I136").Value = 2500
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 7500
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -4950
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1573.8462
$ws.Range("I94").Value = 1192.8572
$ws.Range("K94").Value = 1192.8572
$ws.Range("M94").Value = -741.8571999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1925.7142
$ws.Range("I16").Value = 1896
$ws.Range("K16").Value = 1896
$ws.Range("M16").Value = -1609
# Row 43
$ws.Range("H43").Value = 23749.5
$ws.Range("I43").Value = 10000
$ws.Range("J43").Value = 28332.666
$ws.Range("K43").Value = 10000
$ws.Range("L43").Value = 28332.666
$ws.Range("M43").Value = -9816
$ws.Range("N43").Value = -28700.666
# Row 62
$ws.Range("H62").Value = 102682
$ws.Range("I62").Value = 102682
$ws.Range("K62").Value = 102682
$ws.Range("M62").Value = -102058
# Row 65
$ws.Range("H65").Value = 102682
$ws.Range("I65").Value = 102682
$ws.Range("K65").Value = 513410
$ws.Range("M65").Value = -510290
# Row 101
$ws.Range("H101").Value = 23749.5
$ws.Range("I101").Value = 10000
$ws.Range("J101").Value = 28332.666
$ws.Range("K101").Value = 10000
$ws.Range("L101").Value = 28332.666
$ws.Range("M101").Value = -6755
$ws.Range("N101").Value = -34822.666
# Row 103
$ws.Range("H103").Value = 12753.8
$ws.Range("I103").Value = 1256.3334
$ws.Range("J103").Value = 30000
$ws.Range("K103").Value = 1256.3334
$ws.Range("L103").Value = 30000
$ws.Range("M103").Value = -84.33339999999998
$ws.Range("N103").Value = -32344
# Row 107
$ws.Range("H107").Value = 330.54285
$ws.Range("I107").Value = 282.04544
$ws.Range("J107").Value = 412.6154
$ws.Range("K107").Value = 282.04544
$ws.Range("L107").Value = 412.6154
$ws.Range("M107").Value = 1637.95456
$ws.Range("N107").Value = -4252.6154
# Row 113
$ws.Range("H113").Value = 1925.7142
$ws.Range("I113").Value = 1896
$ws.Range("K113").Value = 1896
$ws.Range("M113").Value = 274
# Row 134
$ws.Range("H134").Value = 1049.7028
$ws.Range("I134").Value = 965.85297
$ws.Range("K134").Value = 2897.55891
$ws.Range("M134").Value = -362.5589100000002

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 517.1905
$ws.Range("I107").Value = 734.2857
$ws.Range("J107").Value = 408.64285
$ws.Range("K107").Value = 2202.8571
$ws.Range("L107").Value = 1225.92855
$ws.Range("M107").Value = -282.8571000000002
$ws.Range("N107").Value = -5065.928550000001
# Row 131
$ws.Range("H131").Value = 16131956
$ws.Range("J131").Value = 18183250
$ws.Range("L131").Value = 54549750
$ws.Range("N131").Value = -54559830
# Row 132
$ws.Range("H132").Value = 1427.0358
$ws.Range("I132").Value = 933.75
$ws.Range("K132").Value = 8403.75
$ws.Range("M132").Value = -5873.75
# Row 137
$ws.Range("H137").Value = 27779726
$ws.Range("J137").Value = 47621776
$ws.Range("L137").Value = 142865328
$ws.Range("N137").Value = -142875528

$ws = $wb.Worksheets.Item("GSM")
# Row 105
$ws.Range("H105").Value = 48450
$ws.Range("J105").Value = 48450
$ws.Range("L105").Value = 48450
$ws.Range("N105").Value = -55438
# Row 122
$ws.Range("H122").Value = 2658.7778
$ws.Range("I122").Value = 2757.261
$ws.Range("K122").Value = 8271.782999999999
$ws.Range("M122").Value = -5821.782999999999
# Row 132
$ws.Range("H132").Value = 3628
$ws.Range("I132").Value = 2942.3333
$ws.Range("K132").Value = 8826.999899999999
$ws.Range("M132").Value = -6296.999899999999
# Row 136
$ws.Range("H136").Value = 13415.542
$ws.Range("J136").Value = 13415.542
$ws.Range("L136").Value = 40246.626
$ws.Range("N136").Value = -45346.626

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3809.5293
$ws.Range("I7").Value = 2823.8572
$ws.Range("J7").Value = 4499.5
$ws.Range("K7").Value = 2823.8572
$ws.Range("L7").Value = 4499.5
$ws.Range("M7").Value = -2711.8572
$ws.Range("N7").Value = -4723.5
# Row 16
$ws.Range("H16").Value = 6950
$ws.Range("I16").Value = 6950
$ws.Range("K16").Value = 6950
$ws.Range("M16").Value = -6780
# Row 93
$ws.Range("H93").Value = 1431.6666
$ws.Range("I93").Value = 1250
$ws.Range("J93").Value = 1795
$ws.Range("K93").Value = 1250
$ws.Range("L93").Value = 1795
$ws.Range("M93").Value = -2
$ws.Range("N93").Value = -4291
# Row 126
$ws.Range("H126").Value = 3809.5293
$ws.Range("I126").Value = 2823.8572
$ws.Range("J126").Value = 4499.5
$ws.Range("K126").Value = 8471.571599999999
$ws.Range("L126").Value = 13498.5
$ws.Range("M126").Value = -6001.571599999999
$ws.Range("N126").Value = -18438.5

$ws = $wb.Worksheets.Item("WVR")
# Row 49
$ws.Range("H49").Value = 14411.2
$ws.Range("J49").Value = 14500
$ws.Range("L49").Value = 14500
$ws.Range("N49").Value = -14960
# Row 54
$ws.Range("H54").Value = 30745
$ws.Range("J54").Value = 30745
$ws.Range("L54").Value = 30745
$ws.Range("N54").Value = -31785
# Row 113
$ws.Range("H113").Value = 857.2
$ws.Range("I113").Value = 659.2857
$ws.Range("J113").Value = 1109.091
$ws.Range("K113").Value = 1977.8571
$ws.Range("L113").Value = 3327.273
$ws.Range("M113").Value = 192.1428999999998
$ws.Range("N113").Value = -7667.272999999999
# Row 122
$ws.Range("H122").Value = 9262516
$ws.Range("I122").Value = 13160214
$ws.Range("K122").Value = 39480642
$ws.Range("M122").Value = -39478192
# Row 126
$ws.Range("H126").Value = 7557.6665
$ws.Range("I126").Value = 9610.223
$ws.Range("J126").Value = 1400
$ws.Range("K126").Value = 28830.669
$ws.Range("L126").Value = 4200
$ws.Range("M126").Value = -26360.669
$ws.Range("N126").Value = -9140
# Row 132
$ws.Range("H132").Value = 2107.889
$ws.Range("I132").Value = 1495
$ws.Range("J132").Value = 3333.6667
$ws.Range("K132").Value = 4485
$ws.Range("L132").Value = 10001.0001
$ws.Range("M132").Value = -1955
$ws.Range("N132").Value = -15061.0001
